# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# The "Periodo Mora" (column E) list of periods is refreshed: the previous
# run of periods (2312, 2401-2406, oldest-to-newest) is replaced with the
# new run of periods in newest-to-oldest order (2406, 2405, 2404, 2403,
# 2402, 2401, 2312), and the "Valor Mora" (column F) figure that belonged
# to period 2406 (18560) now travels with it to row 16, while the rest of
# the rows keep the standard 46400 value (including the row that now shows
# period 2312).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$periods = @("2406", "2405", "2404", "2403", "2402", "2401", "2312")
$valores = @(18560, 46400, 46400, 46400, 46400, 46400, 46400)

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periods[$i]
    $ws.Range("F$row").Value = $valores[$i]
}
